$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("7:7").Insert()

$ws.Range("A7").Value = 1
$ws.Range("B7").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C7").Value = "Arica y Parinacota"
$ws.Range("D7").Value = 44742
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = 100112013
$ws.Range("G7").Value = "Alcachofa"
$ws.Range("H7").Value = "Madrigal"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 120
$ws.Range("K7").Value = 19000
$ws.Range("L7").Value = 20000
$ws.Range("M7").Value = 19500
$ws.Range("N7").Value = '$/caja 40 unidades'
$ws.Range("O7").Value = "Región de Coquimbo"
$ws.Range("P7").Value = 488
$ws.Range("Q7").Value = 40
$ws.Range("R7").Value = "Hortaliza"
